$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2025/12/03 01:25"
$ws.Range("B5").Value = "36,668位本"
$ws.Range("C5").Value = "84位 広告・宣伝 (本)"
$ws.Range("D5").Value = "165位商業デザイン"
$ws.Range("E5").Value = "1,838位ビジネス実用本"
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "-"
